$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (3501 North Lancaster Hutchins Road / Lancaster / TX / 75134)
$ws.Rows.Item(1).Delete()

# After the shift, delete rows 5 through 9 (Tampa, Clayton, Henderson, Aurora, Baton Rouge)
$ws.Range("A5:XFD9").Delete()

# Update selection to match the final state
$ws.Range("A5:XFD9").Select()
